$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 73

# Column A holds a period label formatted like a date ("01-07-2021").
# Force it to be stored as text (matching the other period cells) rather
# than letting Excel auto-convert it to a date serial number.
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "01-07-2021"
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 59.4
$ws.Cells.Item($row, 3).Value = 41.9
$ws.Cells.Item($row, 4).Value = 60.2
$ws.Cells.Item($row, 5).Value = 72.3
$ws.Cells.Item($row, 6).Value = 57
$ws.Cells.Item($row, 7).Value = 67.8
$ws.Cells.Item($row, 8).Value = 73.09999999999999
$ws.Cells.Item($row, 9).Value = 68.59999999999999
$ws.Cells.Item($row, 10).Value = 84.59999999999999
$ws.Cells.Item($row, 11).Value = 58.8
$ws.Cells.Item($row, 12).Value = 63.3
$ws.Cells.Item($row, 13).Value = 61.9
$ws.Cells.Item($row, 14).Value = 62
$ws.Cells.Item($row, 15).Value = 4
$ws.Cells.Item($row, 16).Value = 62.1
